$wb = $excel.ActiveWorkbook

# Update sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9331
$ws1.Range("F4").Value = 19
$ws1.Range("F5").Value = 507

# Update sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9331
$ws4.Range("F4").Value = 19
$ws4.Range("F5").Value = 507
